$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style/number format/borders) from the last existing data row (269)
# down across the new rows 270:301 for columns A:D, matching the established pattern.
$ws.Range("A269:D269").Copy() | Out-Null
$ws.Range("A270:D301").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Write the new daily values for 28/05/2021 (serial 44344) through 28/06/2021 (serial 44375).
$ws.Cells.Item(270,1).Value = 44344
$ws.Cells.Item(270,2).Value = 0
$ws.Cells.Item(270,3).Value = 0
$ws.Cells.Item(270,4).Value = 0
$ws.Cells.Item(271,1).Value = 44345
$ws.Cells.Item(271,2).Value = 0
$ws.Cells.Item(271,3).Value = 0
$ws.Cells.Item(271,4).Value = 0
$ws.Cells.Item(272,1).Value = 44346
$ws.Cells.Item(272,2).Value = 0
$ws.Cells.Item(272,3).Value = 0
$ws.Cells.Item(272,4).Value = 0
$ws.Cells.Item(273,1).Value = 44347
$ws.Cells.Item(273,2).Value = 0
$ws.Cells.Item(273,3).Value = 0
$ws.Cells.Item(273,4).Value = 0
$ws.Cells.Item(274,1).Value = 44348
$ws.Cells.Item(274,2).Value = 0
$ws.Cells.Item(274,3).Value = 0
$ws.Cells.Item(274,4).Value = 0
$ws.Cells.Item(275,1).Value = 44349
$ws.Cells.Item(275,2).Value = 0
$ws.Cells.Item(275,3).Value = 0
$ws.Cells.Item(275,4).Value = 0
$ws.Cells.Item(276,1).Value = 44350
$ws.Cells.Item(276,2).Value = 0
$ws.Cells.Item(276,3).Value = 0
$ws.Cells.Item(276,4).Value = 0
$ws.Cells.Item(277,1).Value = 44351
$ws.Cells.Item(277,2).Value = 0
$ws.Cells.Item(277,3).Value = 0
$ws.Cells.Item(277,4).Value = 0
$ws.Cells.Item(278,1).Value = 44352
$ws.Cells.Item(278,2).Value = 0
$ws.Cells.Item(278,3).Value = 0
$ws.Cells.Item(278,4).Value = 0
$ws.Cells.Item(279,1).Value = 44353
$ws.Cells.Item(279,2).Value = 0
$ws.Cells.Item(279,3).Value = 0
$ws.Cells.Item(279,4).Value = 0
$ws.Cells.Item(280,1).Value = 44354
$ws.Cells.Item(280,2).Value = 0
$ws.Cells.Item(280,3).Value = 0
$ws.Cells.Item(280,4).Value = 0
$ws.Cells.Item(281,1).Value = 44355
$ws.Cells.Item(281,2).Value = 0
$ws.Cells.Item(281,3).Value = 0
$ws.Cells.Item(281,4).Value = 0
$ws.Cells.Item(282,1).Value = 44356
$ws.Cells.Item(282,2).Value = 0
$ws.Cells.Item(282,3).Value = 0
$ws.Cells.Item(282,4).Value = 0
$ws.Cells.Item(283,1).Value = 44357
$ws.Cells.Item(283,2).Value = 0
$ws.Cells.Item(283,3).Value = 0
$ws.Cells.Item(283,4).Value = 0
$ws.Cells.Item(284,1).Value = 44358
$ws.Cells.Item(284,2).Value = 0
$ws.Cells.Item(284,3).Value = 0
$ws.Cells.Item(284,4).Value = 0
$ws.Cells.Item(285,1).Value = 44359
$ws.Cells.Item(285,2).Value = 0
$ws.Cells.Item(285,3).Value = 0
$ws.Cells.Item(285,4).Value = 0
$ws.Cells.Item(286,1).Value = 44360
$ws.Cells.Item(286,2).Value = 1
$ws.Cells.Item(286,3).Value = 1
$ws.Cells.Item(286,4).Value = 33.71544167228591
$ws.Cells.Item(287,1).Value = 44361
$ws.Cells.Item(287,2).Value = 0
$ws.Cells.Item(287,3).Value = 1
$ws.Cells.Item(287,4).Value = 33.71544167228591
$ws.Cells.Item(288,1).Value = 44362
$ws.Cells.Item(288,2).Value = 2
$ws.Cells.Item(288,3).Value = 3
$ws.Cells.Item(288,4).Value = 101.1463250168577
$ws.Cells.Item(289,1).Value = 44363
$ws.Cells.Item(289,2).Value = 0
$ws.Cells.Item(289,3).Value = 3
$ws.Cells.Item(289,4).Value = 101.1463250168577
$ws.Cells.Item(290,1).Value = 44364
$ws.Cells.Item(290,2).Value = 0
$ws.Cells.Item(290,3).Value = 3
$ws.Cells.Item(290,4).Value = 101.1463250168577
$ws.Cells.Item(291,1).Value = 44365
$ws.Cells.Item(291,2).Value = 2
$ws.Cells.Item(291,3).Value = 5
$ws.Cells.Item(291,4).Value = 168.5772083614295
$ws.Cells.Item(292,1).Value = 44366
$ws.Cells.Item(292,2).Value = 0
$ws.Cells.Item(292,3).Value = 5
$ws.Cells.Item(292,4).Value = 168.5772083614295
$ws.Cells.Item(293,1).Value = 44367
$ws.Cells.Item(293,2).Value = 0
$ws.Cells.Item(293,3).Value = 4
$ws.Cells.Item(293,4).Value = 134.8617666891436
$ws.Cells.Item(294,1).Value = 44368
$ws.Cells.Item(294,2).Value = 1
$ws.Cells.Item(294,3).Value = 5
$ws.Cells.Item(294,4).Value = 168.5772083614295
$ws.Cells.Item(295,1).Value = 44369
$ws.Cells.Item(295,2).Value = 0
$ws.Cells.Item(295,3).Value = 3
$ws.Cells.Item(295,4).Value = 101.1463250168577
$ws.Cells.Item(296,1).Value = 44370
$ws.Cells.Item(296,2).Value = 0
$ws.Cells.Item(296,3).Value = 3
$ws.Cells.Item(296,4).Value = 101.1463250168577
$ws.Cells.Item(297,1).Value = 44371
$ws.Cells.Item(297,2).Value = 0
$ws.Cells.Item(297,3).Value = 3
$ws.Cells.Item(297,4).Value = 101.1463250168577
$ws.Cells.Item(298,1).Value = 44372
$ws.Cells.Item(298,2).Value = 0
$ws.Cells.Item(298,3).Value = 1
$ws.Cells.Item(298,4).Value = 33.71544167228591
$ws.Cells.Item(299,1).Value = 44373
$ws.Cells.Item(299,2).Value = 0
$ws.Cells.Item(299,3).Value = 1
$ws.Cells.Item(299,4).Value = 33.71544167228591
$ws.Cells.Item(300,1).Value = 44374
$ws.Cells.Item(300,2).Value = 0
$ws.Cells.Item(300,3).Value = 1
$ws.Cells.Item(300,4).Value = 33.71544167228591
$ws.Cells.Item(301,1).Value = 44375
$ws.Cells.Item(301,2).Value = 0
$ws.Cells.Item(301,3).Value = 0
$ws.Cells.Item(301,4).Value = 0
